$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the data values in row 5 (B5:AH5) to 2 decimal places.
# We use Excel's own ROUND() worksheet function (via a scratch formula
# cell) so the resulting IEEE-754 representation matches exactly what
# Excel itself produces when rounding these numbers.
$scratchRow = 100
for ($col = 2; $col -le 34; $col++) {
    $srcCell = $ws.Cells.Item(5, $col)
    $scratchCell = $ws.Cells.Item($scratchRow, $col)
    $addr = $srcCell.Address($false, $false)
    $scratchCell.Formula = "=ROUND(" + $addr + ",2)"
}
for ($col = 2; $col -le 34; $col++) {
    $scratchCell = $ws.Cells.Item($scratchRow, $col)
    $rounded = $scratchCell.Value()
    $srcCell = $ws.Cells.Item(5, $col)
    $srcCell.Value = $rounded
    $scratchCell.ClearContents()
}

# The last data row (row 6) is removed entirely, shrinking the sheet's
# used range from A1:AH6 down to A1:AH5.
$ws.Rows.Item(6).Delete()
